$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.343.39"
$ws.Range("E2").Value = "  -3.87%  "
$ws.Range("D3").Value = "1.666.10"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.72"
$ws.Range("E5").Value = "  -2.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5163"
$ws.Range("E6").Value = "  -3.33%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06453"
$ws.Range("E8").Value = "  -2.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2566"
$ws.Range("E9").Value = "  -3.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.92"
$ws.Range("E10").Value = "  -4.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07660"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "1.674.56"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").Value = "1.895.35"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.311"
$ws.Range("E14").Value = "  -5.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5546"
$ws.Range("E15").Value = "  -3.90%  "
$ws.Range("D16").Value = "0.0₅8041"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.49"
$ws.Range("E17").Value = "  -5.12%  "
$ws.Range("D18").Value = "26.368.53"
$ws.Range("E18").Value = "  -3.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.33"
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.403"
$ws.Range("E21").Value = "  -5.90%  "
$ws.Range("E22").Value = "  -3.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.902"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.009"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.60"
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.756"
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1160"
$ws.Range("E27").Value = "  -4.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.992"
$ws.Range("E28").Value = "  -4.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.76"
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.378"
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.224"
$ws.Range("E33").Value = "  -6.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.564"
$ws.Range("E34").Value = "  -5.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.759"
$ws.Range("E35").Value = "  -4.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.379"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9252"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5721"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").Value = "1.157.27"
$ws.Range("E39").Value = "  +10.75%  "
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.008"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8438"
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.652"
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.91"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").Value = "1.805.50"
$ws.Range("E45").Value = "  -2.64%  "
$ws.Range("D46").Value = "0.0₈111"
$ws.Range("E46").Value = "  -7.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4504"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.14"
$ws.Range("E48").Value = "  -3.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.922"
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05103"
$ws.Range("E51").Value = "  -2.62%  "
